$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: trigger "Before Insert On Ligne" was renamed/repurposed to
# "Before Insert On PropositionJoueur", assigned to Pierrick
$ws.Range("D5").Value = "Before Insert On PropositionJoueur"
$ws.Range("E5").Value = "Pierrick"

# Row 9: constraint type switched from trigger to procedure, with a new
# (wrapped, multi-line) description instead of a trigger name; the old
# "assigned to" cell (E9) is cleared
$ws.Range("C9").Value = "procédure"
$ws.Range("C9").HorizontalAlignment = -4108  # xlCenter

$ws.Range("D9").Value = "tester nombre emplacements billes" + [char]10 + "et nombre effectif de bille"
$ws.Range("D9").Font.Color = $ws.Range("D5").Font.Color
$ws.Range("D9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 30

$ws.Range("E9").ClearContents()

# Row 10: new "check" mark in column A, and the trigger now fires on Partie
$ws.Range("A10").Value = "c"
$ws.Range("D10").Value = "Before Insert On Partie"

# Row 11: trigger moment switched from Before to After
$ws.Range("D11").Value = "After Insert On Partie"

# New row 16: extra procedure constraint
$ws.Range("B16").Value = "créer une combinaison avec aucune bille en double"
$ws.Range("C16").Value = "procédure"
$ws.Range("C16").NumberFormat = $ws.Range("B3").NumberFormat
$ws.Range("C16").HorizontalAlignment = -4108  # xlCenter

# Column D is now wider to accommodate the longer wrapped text
$ws.Columns.Item(4).ColumnWidth = 37.42578125

# Restore the selection as left by the author
$ws.Range("E11").Select() | Out-Null
